$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price (D) / volume-change (E) columns cell by cell.
# D-column values are plain text in the workbook (e.g. '27.004.86' has
# two dots, and even single-dot-looking values like '1.002' must stay
# TEXT, not be coerced to a number by Excel's smart entry). We force
# text by setting NumberFormat to '@' before assigning the value, then
# reset the cell style back to Normal so no stray style sticks around.

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.004.86'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -3.30%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.714.66'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -3.00%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.03%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '308.12'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -6.25%  '

# Row 6
$ws.Range("E6").Value = '  +0.06%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4742'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +4.15%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3468'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.57%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '42.03'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.07%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07235'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.00%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.040'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -5.13%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.79'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -4.58%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.827'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.98%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '1.715.42'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.94%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.836'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -4.89%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '86.55'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -6.54%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001036'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.37%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06373'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.04%  '

# Row 20
$ws.Range("E20").Value = '  +0.09%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.45'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.84%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.611'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.79%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '27.063.90'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.18%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.72'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -4.35%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.093'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.40%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '19.90'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.16%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '150.67'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -5.88%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.915.24'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.89%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.079'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.89%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '120.42'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.13%  '

# Row 31
$ws.Range("E31").Value = '  -4.96%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09135'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.04%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.602'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.58%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.310'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -5.31%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.469'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +6.49%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02175'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -4.58%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.05864'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -4.26%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.1998'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.52%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '10.94'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -7.59%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.001'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.09%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.711'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -4.62%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5972'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.61%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.083'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -7.94%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '7.459'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -4.47%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '12.77'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.86%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.568'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.37%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5573'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.73%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '118.93'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.27%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.825'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -5.64%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.119'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.11%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06626'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.87%  '

